$wb = $excel.ActiveWorkbook

# --- Sheet: GBV Masibambisane Partners ---
$ws = $wb.Worksheets.Item("GBV Masibambisane Partners")
$ws.Range("C6").Value = "KwaZulu-Natal"
$ws.Range("C7").Value = "KwaZulu-Natal"
$ws.Range("G7").Value = -29.56617368768535
$ws.Range("H7").Value = 30.18582508164073
$ws.Range("C8").Value = "KwaZulu-Natal"
$ws.Range("G8").Value = -29.56526933488195
$ws.Range("H8").Value = 30.1768329958818
$ws.Range("G9").Value = -23.30472086703238
$ws.Range("H9").Value = 30.70747633816771
$ws.Range("G10").Value = -23.30549895907519
$ws.Range("H10").Value = 30.70299779643534
$ws.Range("C14").Value = "Western Cape"
$ws.Range("C23").Value = "Western Cape"

# --- Sheet: Growing Food Partners ---
$ws = $wb.Worksheets.Item("Growing Food Partners")
$ws.Range("C3").Value = "KwaZulu-Natal"
$ws.Range("C11").Value = "North West"
$ws.Range("C14").Value = "KwaZulu-Natal"
$ws.Range("C17").Value = "KwaZulu-Natal"
$ws.Range("C18").Value = "KwaZulu-Natal"
$ws.Range("C20").Value = "North West"

# --- Sheet: Human Rights Clubs ---
$ws = $wb.Worksheets.Item("Human Rights Clubs")
$ws.Range("H2").Value = -25.67561525692701
$ws.Range("I2").Value = 27.23584568333255
$ws.Range("H3").Value = -25.67639563630532
$ws.Range("I3").Value = 27.24686548918428
$ws.Range("H4").Value = -33.91486044871371
$ws.Range("I4").Value = 18.4303868407082
$ws.Range("H5").Value = -33.91843687622078
$ws.Range("I5").Value = 18.42483316221153
$ws.Range("H6").Value = -23.84121903436854
$ws.Range("I6").Value = 29.37694249305345
$ws.Range("H7").Value = -23.84114844892832
$ws.Range("I7").Value = 29.37960026836785
$ws.Range("E8").Value = "Eastern Cape"
$ws.Range("H8").Value = -33.90650348116984
$ws.Range("I8").Value = 25.58779022272385
$ws.Range("E9").Value = "Eastern Cape"
$ws.Range("H9").Value = -33.90675649172453
$ws.Range("I9").Value = 25.58837571349237
$ws.Range("H10").Value = -30.65129098950961
$ws.Range("I10").Value = 24.01401732006992
$ws.Range("H11").Value = -30.63777962344104
$ws.Range("I11").Value = 24.00869877999387
$ws.Range("E12").Value = "KwaZulu-Natal"
$ws.Range("H12").Value = -29.59647196517574
$ws.Range("I12").Value = 30.37921417005645
$ws.Range("E13").Value = "KwaZulu-Natal"
$ws.Range("H13").Value = -29.60054737026591
$ws.Range("I13").Value = 30.37802430166695
$ws.Range("H14").Value = -29.12758616489118
$ws.Range("I14").Value = 26.28562167835697
$ws.Range("H15").Value = -29.13081182836828
$ws.Range("I15").Value = 26.28471596254374
$ws.Range("H16").Value = -29.13979741335006
$ws.Range("I16").Value = 26.28036611426643
$ws.Range("H17").Value = -29.13312496661504
$ws.Range("I17").Value = 26.28141291430465
$ws.Range("H18").Value = -26.01617497365869
$ws.Range("I18").Value = 28.22529983485829
$ws.Range("H19").Value = -26.00772897785614
$ws.Range("I19").Value = 28.21686510370403
$ws.Range("H20").Value = -26.35536144556415
$ws.Range("I20").Value = 28.1443528409034
$ws.Range("H21").Value = -26.36498227083928
$ws.Range("I21").Value = 28.14407618099457
$ws.Range("H22").Value = -25.47635795003008
$ws.Range("I22").Value = 29.00746233913173
$ws.Range("H23").Value = -25.4863984750171
$ws.Range("I23").Value = 29.02125058596489
$ws.Range("H24").Value = -31.59129237924279
$ws.Range("I24").Value = 28.76622748996408
$ws.Range("H25").Value = -31.59383582131341
$ws.Range("I25").Value = 28.75713086575806
$ws.Range("H26").Value = -23.30699409606289
$ws.Range("I26").Value = 30.69488038472128
$ws.Range("H27").Value = -23.31121455460828
$ws.Range("I27").Value = 30.6948092315841
$ws.Range("H28").Value = -25.48663210322576
$ws.Range("I28").Value = 27.83662532363554
$ws.Range("H29").Value = -25.48886711686297
$ws.Range("I29").Value = 27.83626665444433
$ws.Range("E30").Value = "KwaZulu-Natal"
$ws.Range("H30").Value = -29.55937289129131
$ws.Range("I30").Value = 30.19448471978172
$ws.Range("E31").Value = "KwaZulu-Natal"
$ws.Range("H31").Value = -29.57737026179612
$ws.Range("I31").Value = 30.1887879545572
$ws.Range("H32").Value = -26.03136101452483
$ws.Range("I32").Value = 30.785120257149
$ws.Range("H33").Value = -26.04971293877983
$ws.Range("I33").Value = 30.79186298832423
$ws.Range("H34").Value = -24.8343291472963
$ws.Range("I34").Value = 31.07503226936309
$ws.Range("H35").Value = -24.83413127169455
$ws.Range("I35").Value = 31.07595594527941
$ws.Range("E36").Value = "Eastern Cape"
$ws.Range("H36").Value = -31.91989986631074
$ws.Range("I36").Value = 26.96835760417916
$ws.Range("H37").Value = -31.92440641570071
$ws.Range("I37").Value = 26.96170478174716
$ws.Range("E38").Value = "Eastern Cape"
$ws.Range("H38").Value = -31.91903281002792
$ws.Range("I38").Value = 26.95859839517394
$ws.Range("H39").Value = -29.22807044592622
$ws.Range("I39").Value = 26.7137045122598
$ws.Range("H40").Value = -29.23409597626011
$ws.Range("I40").Value = 26.7050172654635
$ws.Range("H41").Value = -28.73917847099753
$ws.Range("I41").Value = 24.76517568009713
$ws.Range("H42").Value = -28.74098160667427
$ws.Range("I42").Value = 24.76734260845445
$ws.Range("H43").Value = -33.93523786602753
$ws.Range("I43").Value = 18.85313351420969
$ws.Range("H44").Value = -33.92743585743257
$ws.Range("I44").Value = 18.86326339619207
$ws.Range("E45").Value = "KwaZulu-Natal"
$ws.Range("H45").Value = -29.57466070029357
$ws.Range("I45").Value = 30.18415094608206
$ws.Range("E46").Value = "KwaZulu-Natal"
$ws.Range("H46").Value = -29.56073443345731
$ws.Range("I46").Value = 30.17832181777939
$ws.Range("H47").Value = -23.32126218644981
$ws.Range("I47").Value = 30.70499327015509
$ws.Range("H48").Value = -23.30313154721219
$ws.Range("I48").Value = 30.70852188236656

# --- Sheet: TRC Cases Supported ---
$ws = $wb.Worksheets.Item("TRC Cases Supported")
$ws.Range("C7").Value = "North West"
$ws.Range("H9").Value = -26.20779877670934
$ws.Range("I9").Value = 28.03384461945495
$ws.Range("H12").Value = -26.2021863746939
$ws.Range("I12").Value = 28.03158767278949
$ws.Range("H13").Value = -29.86424184272459
$ws.Range("I13").Value = 31.0229213815519
$ws.Range("H14").Value = -33.92560180897983
$ws.Range("I14").Value = 18.41867477926927
$ws.Range("H15").Value = -26.19961059762062
$ws.Range("I15").Value = 28.02735632037809
$ws.Range("H18").Value = -33.97358883456565
$ws.Range("I18").Value = 18.50096177489309
$ws.Range("H19").Value = -29.85373171138816
$ws.Range("I19").Value = 31.01856758092833
$ws.Range("H20").Value = -33.92638245895357
$ws.Range("I20").Value = 18.42920592122096
$ws.Range("H21").Value = -33.96541170497592
$ws.Range("I21").Value = 18.50170394113072
$ws.Range("H22").Value = -26.20219738747916
$ws.Range("I22").Value = 28.03262283151156
$ws.Range("H24").Value = -29.86484746922991
$ws.Range("I24").Value = 31.02622780455514
